$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2
$ws1.Range("H2").Value = 16.06
$ws1.Range("L2").Value = 1.1

# Row 3
$ws1.Range("H3").Value = 12.24
$ws1.Range("L3").Value = 0.98

# Row 4
$ws1.Range("H4").Value = 13.83
$ws1.Range("L4").Value = 0.99

# Row 5
$ws1.Range("H5").Value = 10.43
$ws1.Range("L5").Value = 0.9

# Row 6
$ws1.Range("H6").Value = 9.43
$ws1.Range("L6").Value = 0.91

# Row 7
$ws1.Range("H7").Value = 8.43
$ws1.Range("L7").Value = 1.06

# Row 8
$ws1.Range("H8").Value = 7.43
$ws1.Range("L8").Value = 0.83

# Row 9
$ws1.Range("H9").Value = 6.43
$ws1.Range("L9").Value = 0.97

# Row 10
$ws1.Range("D10").Value = 4
$ws1.Range("H10").Value = 5.43
$ws1.Range("L10").Value = 1.01

# Row 11
$ws1.Range("D11").Value = 4
$ws1.Range("H11").Value = 4.43
$ws1.Range("L11").Value = 0.9

# Row 12
$ws1.Range("D12").Value = 4
$ws1.Range("H12").Value = 3.43
$ws1.Range("L12").Value = 0.9

# Row 13
$ws1.Range("H13").Value = 2.43
$ws1.Range("L13").Value = 0.97

# Row 14
$ws1.Range("D14").Value = 3
$ws1.Range("H14").Value = 1.57
$ws1.Range("L14").Value = 0.95

# Row 15
$ws1.Range("H15").Value = 0.57
$ws1.Range("I15").Value = "Low"
$ws1.Range("L15").Value = 1.19

# Row 16
$ws1.Range("L16").Value = 1.08

# Row 17
$ws1.Range("L17").Value = 1.02

# --- Sheet: Summary ---
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "62"
$ws2.Range("B10").Value = "31"
$ws2.Range("B11").Value = "15"
$ws2.Range("B14").Value = "3"
